$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The file listing should only keep the header row plus the 2012 entry.
# Delete the rows for 2002-2003 and 2004-2023 (rows 3 through 23), which
# shifts everything up and shrinks the used range to A1:B2.
$ws.Range("A3:B23").EntireRow.Delete() | Out-Null

# Row 2 now represents the 2012 data file (previously row 12's values).
$ws.Range("A2").Value = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\2012.xlsx"
$ws.Range("B2").Value = "2012"

# Match the target selection: a single active cell instead of the old A2:B23 block.
$ws.Range("A2").Select()
